# ---------------------------------------------------------------------------
# Scheduled Sheets data refresh
#
# This workbook tracks FFXIV Leve "turn-in" profitability: for every Leve row,
# columns H:N hold market-board derived figures (currentAveragePrice[NQ/HQ],
# LevePrice[NQ/HQ] and LeveProfit[NQ/HQ]). A scheduled runner periodically
# re-pulls current market prices and rewrites those value cells per item/sheet.
# A handful of rows whose HQ price column fell to 0 no longer have a valid
# profit figure, so the corresponding profit cell is cleared outright instead
# of being rewritten.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ===== Sheet "ALC" =====
$ws = $wb.Worksheets.Item("ALC")

# Row 12
$ws.Range("H12").Value = 554.5714
$ws.Range("I12").Value = 497.16666
$ws.Range("J12").Value = 899
$ws.Range("K12").Value = 497.16666
$ws.Range("L12").Value = 899
$ws.Range("M12").Value = -327.16666
$ws.Range("N12").Value = -1239

# Row 39
$ws.Range("H39").Value = 1789.75
$ws.Range("I39").Value = 121
$ws.Range("K39").Value = 363
$ws.Range("M39").Value = -67

# Row 51
$ws.Range("H51").Value = 43679.41
$ws.Range("I51").Value = 10994.556
$ws.Range("J51").Value = 66307.38
$ws.Range("K51").Value = 10994.556
$ws.Range("L51").Value = 66307.38
$ws.Range("M51").Value = -10510.556
$ws.Range("N51").Value = -67275.38

# Row 138
$ws.Range("H138").Value = 2286.5151
$ws.Range("J138").Value = 3685.8333
$ws.Range("L138").Value = 11057.4999
$ws.Range("N138").Value = -21337.4999

# ===== Sheet "ARM" =====
$ws = $wb.Worksheets.Item("ARM")

# Row 5
$ws.Range("H5").Value = 469
$ws.Range("I5").Value = 458.66666
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 458.66666
$ws.Range("L5").Value = 500
$ws.Range("M5").Value = -346.66666
$ws.Range("N5").Value = -724

# Row 44
$ws.Range("H44").Value = 59618.125
$ws.Range("J44").Value = 59618.125
$ws.Range("L44").Value = 59618.125
$ws.Range("N44").Value = -60594.125

# Row 45
$ws.Range("H45").Value = 1653.5883
$ws.Range("I45").Value = 1582.4
$ws.Range("K45").Value = 1582.4
$ws.Range("M45").Value = -1205.4

# Row 92
$ws.Range("H92").Value = 40909.8
$ws.Range("J92").Value = 40909.8
$ws.Range("L92").Value = 40909.8
$ws.Range("N92").Value = -45901.8

# ===== Sheet "BSM" =====
$ws = $wb.Worksheets.Item("BSM")

# Row 4
$ws.Range("H4").Value = 469
$ws.Range("I4").Value = 458.66666
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 458.66666
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = -343.66666
$ws.Range("N4").Value = -730

# Row 15
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()

# Row 19
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()

# Row 20
$ws.Range("H20").Value = 2656.5
$ws.Range("I20").Value = 2522.5715
$ws.Range("J20").Value = 2890.875
$ws.Range("K20").Value = 2522.5715
$ws.Range("L20").Value = 2890.875
$ws.Range("M20").Value = -2275.5715
$ws.Range("N20").Value = -3384.875

# Row 105
$ws.Range("H105").Value = 50013804
$ws.Range("I105").Value = 66684100
$ws.Range("K105").Value = 66684100
$ws.Range("M105").Value = -66682353

# ===== Sheet "CRP" =====
$ws = $wb.Worksheets.Item("CRP")

# Row 58
$ws.Range("H58").Value = 5549.727
$ws.Range("I58").Value = 4716
$ws.Range("K58").Value = 4716
$ws.Range("M58").Value = -4513

# Row 62
$ws.Range("H62").Value = 6828.5
$ws.Range("I62").Value = 3914
$ws.Range("J62").Value = 9743
$ws.Range("K62").Value = 3914
$ws.Range("L62").Value = 9743
$ws.Range("M62").Value = -3290
$ws.Range("N62").Value = -10991

# Row 65
$ws.Range("H65").Value = 6828.5
$ws.Range("I65").Value = 3914
$ws.Range("J65").Value = 9743
$ws.Range("K65").Value = 19570
$ws.Range("L65").Value = 48715
$ws.Range("M65").Value = -16450
$ws.Range("N65").Value = -54955

# Row 133
$ws.Range("H133").Value = 102777.75
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 102777.75
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 102777.75
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -107837.75

# Row 134
$ws.Range("H134").Value = 9804.916999999999
$ws.Range("I134").Value = 10011.875
$ws.Range("K134").Value = 30035.625
$ws.Range("M134").Value = -27500.625

# Row 136
$ws.Range("H136").Value = 5549.727
$ws.Range("I136").Value = 4716
$ws.Range("K136").Value = 14148
$ws.Range("M136").Value = -11598

# ===== Sheet "CUL" =====
$ws = $wb.Worksheets.Item("CUL")

# Row 26
$ws.Range("H26").Value = 7063.125

# Row 58
$ws.Range("H58").Value = 1912.4
$ws.Range("I58").Value = 1021.3333
$ws.Range("J58").Value = 3249
$ws.Range("K58").Value = 3063.9999
$ws.Range("L58").Value = 9747
$ws.Range("M58").Value = -2935.9999
$ws.Range("N58").Value = -10003

# Row 113
$ws.Range("H113").Value = 2361.6875
$ws.Range("J113").Value = 2101.25
$ws.Range("L113").Value = 6303.75
$ws.Range("N113").Value = -10643.75

# ===== Sheet "GSM" =====
$ws = $wb.Worksheets.Item("GSM")

# Row 70
$ws.Range("H70").Value = 3809.0833
$ws.Range("I70").Value = 3809.0833
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 3809.0833
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -3539.0833
$ws.Range("N70").ClearContents()

# Row 73
$ws.Range("H73").Value = 3809.0833
$ws.Range("I73").Value = 3809.0833
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 3809.0833
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -2873.0833
$ws.Range("N73").ClearContents()

# ===== Sheet "LTW" =====
$ws = $wb.Worksheets.Item("LTW")

# Row 22
$ws.Range("H22").Value = 3139.261
$ws.Range("I22").Value = 2258.889
$ws.Range("J22").Value = 3705.2144
$ws.Range("K22").Value = 2258.889
$ws.Range("L22").Value = 3705.2144
$ws.Range("M22").Value = -1963.889
$ws.Range("N22").Value = -4295.2144

# Row 27
$ws.Range("H27").Value = 3139.261
$ws.Range("I27").Value = 2258.889
$ws.Range("J27").Value = 3705.2144
$ws.Range("K27").Value = 2258.889
$ws.Range("L27").Value = 3705.2144
$ws.Range("M27").Value = -2151.889
$ws.Range("N27").Value = -3919.2144

# Row 46
$ws.Range("H46").Value = 7201.222
$ws.Range("J46").Value = 7637.28
$ws.Range("L46").Value = 7637.28
$ws.Range("N46").Value = -8013.28

# Row 55
$ws.Range("I55").Value = 238.6
$ws.Range("J55").Value = 216.16667
$ws.Range("K55").Value = 238.6
$ws.Range("L55").Value = 216.16667
$ws.Range("M55").Value = -65.59999999999999
$ws.Range("N55").Value = -562.1666700000001

# Row 100
$ws.Range("H100").Value = 13161698
$ws.Range("I100").Value = 41669784
$ws.Range("J100").Value = 4120.154
$ws.Range("K100").Value = 41669784
$ws.Range("L100").Value = 4120.154
$ws.Range("M100").Value = -41669243

# Row 139
$ws.Range("H139").Value = 132123.62
$ws.Range("J139").Value = 109498.164
$ws.Range("L139").Value = 109498.164
$ws.Range("N139").Value = -119778.164

# ===== Sheet "WVR" =====
$ws = $wb.Worksheets.Item("WVR")

# Row 132
$ws.Range("H132").Value = 5064
$ws.Range("I132").Value = 3428.6924
$ws.Range("K132").Value = 10286.0772
$ws.Range("M132").Value = -7756.0772
